$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6 (Generation 0-4): Fitness 7310 -> 7295
$ws.Range("C2:C6").Value = 7295

# Rows 7-169 (Generation 5-167): Fitness 7310 -> 7293
$ws.Range("C7:C169").Value = 7293
